# Add IPT in adults
# This script inserts a new time-variant parameter row
# ("program_perc_ipt_age15up") into the "time_variants" sheet, just
# below "program_perc_ipt_age5to15" (old row 16) and above
# "program_perc_awareness_raising" (old row 17), shifting all
# subsequent rows down by one. It also adds a couple of extra
# scenario values to the two existing IPT rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("time_variants")

# --- Add new scenario values (scenario_13 / scenario_14) to the two
#     existing IPT rows (program_perc_ipt_age0to5, program_perc_ipt_age5to15)
$ws.Range("BO15").Value = 90
$ws.Range("BP15").Value = 90
$ws.Range("BO16").Value = 90
$ws.Range("BP16").Value = 90

# --- Insert a new (blank) row at 17 for "program_perc_ipt_age15up",
#     pushing the existing rows 17-27 down to 18-28.
$ws.Rows("17:17").Insert()

# Copy cell formatting from row 15 (program_perc_ipt_age0to5), which
# already has the same look as the new row, onto the freshly inserted
# row 17.
$ws.Range("A15").Copy()
$ws.Range("A17").PasteSpecial(-4122)
$ws.Range("B15:D15").Copy()
$ws.Range("B17:D17").PasteSpecial(-4122)
$ws.Range("BB15").Copy()
$ws.Range("BB17").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0

# The inserted row inherited a few stray blank-styled cells from its
# neighbours (row 16/18); remove them completely since the new row only
# carries data in columns A-D, BB and BO.
$ws.Range("BE17:BP17").Clear()

# Populate the new row's content
$ws.Range("A17").Value = "program_perc_ipt_age15up"
$ws.Range("B17").Value = "no"
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = "yes"
$ws.Range("BB17").Value = 0
$ws.Range("BO17").Value = 90

# Restore the selection to reflect the edited row
$ws.Range("A18").Select() | Out-Null
